$d = $word.ActiveDocument

# "Tanveer Salim (" and "R11597879)" were previously two separate runs
# within the same paragraph. Word's Find/Replace matches text across run
# boundaries and rewrites the match as a single run, so replacing the
# combined text with itself merges the two runs into one.
$d.Content.Find.Execute("Tanveer Salim (R11597879)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Tanveer Salim (R11597879)", 2)

# Add Rafael Perales' R# next to his name, same as Tanveer's line.
$d.Content.Find.Execute("Rafael Perales", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Rafael Perales (R#11458275)", 2)
